# Categorical-inference recap table: add a new "categorical" worksheet after
# Sheet1, populate it, restyle it to match the existing table's font, and
# leave the view/selection state the way the author left it (new sheet
# active, old sheet's selection widened to the full table).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1: selection grows from a single cell to the whole table, and it is
# no longer the tab shown when the file re-opens (the new sheet takes over).
$ws1.Range("A1:E5").Select()

# New sheet, inserted immediately after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "categorical"

# Fill the table in the same order it was authored so the shared-string
# table comes out in the same sequence.
$ws2.Range("D2").Value = '$\hat{p}$'
$ws2.Range("E2").Value = '$N(0,1)$'
$ws2.Range("F2").Value = "prop_test()"

$ws2.Range("A1").Value = "Variable(s) of interest"
$ws2.Range("A2").Value = "One 2-level categorical variable"
$ws2.Range("A3").Value = "Two 2-level categorical variables"
$ws2.Range("A4").Value = "Two 2+ level categorical variables"

$ws2.Range("B2").Value = "One proportion z-test"
$ws2.Range("B1").Value = "Hypothesis Test"
$ws2.Range("B3").Value = "Two proportion z-test"

$ws2.Range("C2").Value = '$p$'
$ws2.Range("C3").Value = '$p_1 - p_2$'
$ws2.Range("C1").Value = "Parameter for CI"
$ws2.Range("C4").Value = "NA"

$ws2.Range("B4").Value = "Chi-square test"

$ws2.Range("D3").Value = '$\hat{p}_1 - \hat{p}_2$'
$ws2.Range("D4").Value = '$X^2$'
$ws2.Range("E4").Value = '$\chi^2_{(I-1)(J-1)}$'
$ws2.Range("F4").Value = "chisq_test() or chisq.test()"

$ws2.Range("D1").Value = "Point Estimate"
$ws2.Range("E1").Value = "Reference distribution"
$ws2.Range("F1").Value = "R function"
$ws2.Range("E3").Value = '$N(0,1)$'
$ws2.Range("F3").Value = "prop_test()"

# Row 5 stays blank but picks up the same styling as the rest of the table.
# Match the black-font style already used on Sheet1 (e.g. D4) across the
# whole new table, including the blank trailing row.
$ws2.Range("A1:F5").Font.Color = 0

# The new sheet becomes the active tab, with F4 as the last-selected cell.
$ws2.Activate()
$ws2.Range("F4").Select()

$wb.Save()
